# Swap the "<<1a1>>" / "<<1a3" placeholder labels inside the "Co quan:"
# run, then re-anchor the document's "_GoBack" bookmark (the marker Word
# drops at the location of the most recent edit) to that spot, removing
# it from its old location near "...vao hoi".

$d = $word.ActiveDocument

# --- 1. Swap the placeholder tokens in the "Cong an <<1a1>>, <<1a3" run ---
$rng = $d.Content
$found = $rng.Find.Execute(
    "<<#1a1>>Công an <<1a1>>, <<1a3",  # FindText
    $true,                             # MatchCase
    $false,                            # MatchWholeWord
    $false,                            # MatchWildcards
    $false,                            # MatchSoundsLike
    $false,                            # MatchAllWordForms
    $true,                             # Forward
    1,                                 # Wrap            (wdFindContinue)
    $false,                            # Format
    "<<#1a1>>Công an <<1a3>>, <<1a1",  # ReplaceWith
    1                                  # Replace         (wdReplaceOne)
)

# --- 2. Move "_GoBack" so it sits right after the text we just edited ---
# (Word always collapses "_GoBack" to a zero-length bookmark at the point
# of the last edit; Bookmarks.Add re-targets it if it already exists,
# which also drops it from its previous location automatically.)
if ($found) {
    $editEnd = $rng.End
    $goBackRange = $d.Range($editEnd, $editEnd)
    $d.Bookmarks.Add("_GoBack", $goBackRange)
}
